$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 (Nate) - all matches played by 10/12
$ws.Range("B8").Value = 39
$ws.Range("C8").Value = 25
$ws.Range("D8").Value = 25
$ws.Range("E8").Value = 325
$ws.Range("F8").Value = 15
$ws.Range("G8").Value = 4335
$ws.Range("H8").Value = 47

# Row 9 (Jsad)
$ws.Range("B9").Value = 28
$ws.Range("C9").Value = 11
$ws.Range("D9").Value = 23
$ws.Range("E9").Value = 92
$ws.Range("F9").Value = 13
$ws.Range("G9").Value = 3336
$ws.Range("H9").Value = 47

# Row 12 (Joey)
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = 27
$ws.Range("D12").Value = 46
$ws.Range("E12").Value = 308
$ws.Range("F12").Value = 20
$ws.Range("G12").Value = 5150
$ws.Range("H12").Value = 60

# Row 13 (Dmah)
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 49
$ws.Range("E13").Value = 111
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 1886
$ws.Range("H13").Value = 60

# Update the active selection to G17
$ws.Range("G17").Select()
